$newText = @'
questions = [
    {
        "title": "While working on a large project with a lot of models and objects, you made some objects invisible to prevent Blender from slowing down. Making every individual item visible again can be a long process.  While in the viewport, what shortcut can solve this problem?",
        "ques_type": 2,
        "options": [
            "Shift + H",
            "Shift + S",
            "Alt + H",
            "Ctrl + S (PC)/Cmd + S (Apple)"
        ],
        "score": "Alt + H"
    },
    {
        "title": "You are working on a disco ball mesh. To improve Blender\u2019s performance, you want to lower the polygon density of the sphere in Image 1 so that it resembles the sphere in Image 2.  While in edit mode, which tool should you use to achieve this?",
        "ques_type": 2,
        "options": [
            "Inset Faces",
            "Reduce Poly Count",
            "Un-subdivide",
            "Bridge Faces"
        ],
        "score": "Un-subdivide"
    },
    {
        "title": "You are creating a scene that features multiple babies that function in the same way. You have completed the weight painting of the first baby. Now, in order to avoid having to repeat the task multiple times, you want to copy the weights to the duplicate characters.  How should you achieve this?",
        "ques_type": 2,
        "options": [
            "Select both models &gt enter Weight Paint Mode &gt click Weight &gt select Mirror.",
            "Select both models &gt enter Weight Paint Mode &gt click Weight &gt select Transfer Weights.",
            "Select the original baby model &gt enter Edit Mode &gt click Weight &gt select Invert.",
            "Select the original baby model &gt select Copy &gt select the duplicate baby model &gt select Paste."
        ],
        "score": "Select both models &gt enter Weight Paint Mode &gt click Weight &gt select Transfer Weights."
    },
    {
        "title": "You have created a human face that looks very realistic, and now you want to render the scene from the camera view in Blender. Currently, the scene is badly lit, with the shadows making it difficult to see the character\u2019s details. To improve this situation, you want to set up three types of light in the scene to illuminate the human face.  Which of the following lights should you use?",
        "ques_type": 15,
        "options": [
            "Sun light",
            "Key light",
            "Area light",
            "Spot light",
            "Rim light",
            "Fill light"
        ],
        "score": [
            "Key light",
            "Rim light",
            "Fill light"
        ]
    }
]
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText
$ws.Rows(1).AutoFit()
